# Added Both ON/OFF vendors in Download SF List
# Adds a new "On/Off" column (AH) -> {vendor:on_off_status} merge field
# to the SF List download template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column data ----------------------------------------------------
# Row 1 header cell - bold + centered, same look as the other header cells
# (e.g. AG1 "Status").
$ws.Range("AH1").Value = "On/Off"
$ws.Range("AH1").Font.Bold = $true
$ws.Range("AH1").HorizontalAlignment = -4108   # xlCenter

# Row 2 placeholder cell - plain style, same as the other merge-field row
# cells (e.g. AB2).
$ws.Range("AH2").Value = "{vendor:on_off_status}"

# --- Column sizing -------------------------------------------------------
# Match the new column's width (~20.22 chars, best-fit) as closely as the
# engine's width quantization allows.
$ws.Columns.Item(34).ColumnWidth = 19.3

# --- View state ------------------------------------------------------------
# Scroll the view over to show the new column and select AI7, mirroring the
# window position captured when the column was added in Excel.
$aw = $excel.ActiveWindow
$aw.ScrollRow = 1
$aw.ScrollColumn = 28
$ws.Range("AI7").Select()
